# Adding new middleman server project and updating other clients with new
# ipaddresses for middleman.
#
# This adds a new weekly-status worksheet "10-7-13" (copied from the most
# recent sheet, "9-30-13", so it keeps the same headers/formatting), fills
# in the new week's task rows, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# The most recently added weekly sheet is the last tab - use it as the
# template for the new week so headers, number formats and column widths
# match the existing sheets exactly.
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet.Copy($null, $templateSheet) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "10-7-13"

# Clear out the copied task rows (2013-09-30 data) - keep rows 1-3 (the
# Date/header rows) which are identical across every weekly sheet.
$newSheet.Range("A4:N7").ClearContents()

# New week's "Date:" header value (10/7/13).
$newSheet.Range("B1").Value = 40092

# Row 4: Reading task carried over from the prior week.
$newSheet.Range("A4").Value = "Read Android Documentation on UDP servers and using both send and receive in same app"
$newSheet.Range("B4").Value = 40088
$newSheet.Range("C4").Value = 40088
$newSheet.Range("D4").Value = 1
$newSheet.Range("E4").Value = 1

# Row 5: Code Android send/receive server.
$newSheet.Range("A5").Value = "Code Android send/receive server"
$newSheet.Range("B5").Value = 40089
$newSheet.Range("C5").Value = 40089
$newSheet.Range("D5").Value = 1
$newSheet.Range("E5").Value = 2.5

# Row 6: Test Android send/receive server.
$newSheet.Range("A6").Value = "Test Android send/receive server"
$newSheet.Range("B6").Value = 40092
$newSheet.Range("C6").Value = 40089
$newSheet.Range("D6").Value = 1
$newSheet.Range("E6").Value = 0.5

# Row 7: Code and test forwarding of movement commands, ignoring others
# for now.
$newSheet.Range("A7").Value = "Code and test forwarding of movement commands, ignoring others for now"
$newSheet.Range("B7").Value = 40102
$newSheet.Range("C7").Value = 40089
$newSheet.Range("D7").Value = 1
$newSheet.Range("E7").Value = 0.5

# Match the selection left on the new sheet (just past the last used row).
$newSheet.Range("A8").Select() | Out-Null

# Make the newly added week the active tab.
$newSheet.Activate() | Out-Null
